{"js": "// almacen.docx \u2014 update the two FSC representative names:\n//   \"Javier Jim\u00e9nez\"  -> \"Betzabet Mar\u00edn\"   (item 3.1, Sistema FSC rep.)\n//   \"Fernando G\u00f3mez\"  -> \"Araceli Becerril\" (item 3.2, Higiene y Seguridad rep.)\n\nconst body = context.document.body;\n\nconst hits1 = body.search(\"Javier Jim\u00e9nez\", { matchCase: true });\nhits1.load(\"items\");\nawait context.sync();\nhits1.items.forEach((r) => r.insertText(\"Betzabet Mar\u00edn\", Word.InsertLocation.replace));\nawait context.sync();\n\nconst hits2 = body.search(\"Fernando G\u00f3mez\", { matchCase: true });\nhits2.load(\"items\");\nawait context.sync();\nhits2.items.forEach((r) => r.insertText(\"Araceli Becerril\", Word.InsertLocation.replace));\nawait context.sync();\n", "ps1": "# almacen.docx \u2014 update the two FSC representative names:\n#   \"Javier Jim\u00e9nez\"  -> \"Betzabet Mar\u00edn\"   (item 3.1, Sistema FSC rep.)\n#   \"Fernando G\u00f3mez\"  -> \"Araceli Becerril\" (item 3.2, Higiene y Seguridad rep.)\n\n$d = $word.ActiveDocument\n\n# --- 1. Replace the Sistema FSC representative's name -----------------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"Javier Jim\u00e9nez\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"Betzabet Mar\u00edn\"\n$find1.Execute([ref]$find1.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$false, [ref]$find1.Replacement.Text, 2)\n\n# --- 2. Replace the Higiene y Seguridad representative's name ---------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Fernando G\u00f3mez\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"Araceli Becerril\"\n$find2.Execute([ref]$find2.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$false, [ref]$find2.Replacement.Text, 2)\n\n# --- 3. Word re-anchors the hidden \"_GoBack\" (last-edit) bookmark to the  ---\n#        new edit site whenever a save happens after text elsewhere in the\n#        document changes. Mirror that bookkeeping: drop the old bookmark\n#        (it used to sit mid-paragraph in item 11, splitting a run in two)\n#        and recreate it collapsed at the start of the paragraph that\n#        follows the text we just edited (\"4-. Debes de conocer...\").\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Delete()\n\n# Re-typing across the old split point lets the two identically-formatted\n# runs in item 11 collapse back into a single run, same as Word does.\n$boundary = $d.Content\n$boundary.Find.ClearFormatting()\n$boundary.Find.Text = \"podemos com\"\n$boundary.Find.Execute() | Out-Null\n$splitPos = $boundary.End\n\n$charRng = $d.Range($splitPos - 1, $splitPos)\n$savedChar = $charRng.Text\n$charRng.Delete()\n$reinsertPt = $d.Range($splitPos - 1, $splitPos - 1)\n$reinsertPt.InsertAfter($savedChar)\n\n$target = $d.Content\n$target.Find.ClearFormatting()\n$target.Find.Text = \"4-. Debes de conocer\"\n$target.Find.Execute() | Out-Null\n$d.Bookmarks.Add(\"_GoBack\", $d.Range($target.Start, $target.Start))\n"}
